# Daily attendance processing - 2025-10-20 03:52:33
# Re-orders the "Recorded By" (column G) contributor list so that any
# entries equal to "system" (case-insensitive) are moved to the front of
# the comma-separated list, while preserving the relative order of the
# "system" entries among themselves and the relative order of the
# remaining (non-system) entries among themselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7 ("Recorded By")
    $val = $cell.Text

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }
    if ($val -notlike "*,*") { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $systemItems = @()
    $otherItems = @()
    foreach ($item in $trimmed) {
        if ($item.ToLower() -eq "system") {
            $systemItems += $item
        } else {
            $otherItems += $item
        }
    }

    $newOrder = $systemItems + $otherItems
    $newVal = [string]::Join(", ", $newOrder)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
